$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Person #4: Jelle Nauta -> Jorrit de Boer
$ws.Range("B6").Value = "Jorrit"
$ws.Range("C6").Value = "de Boer"

# 2. Insert a new row above row 19 so the [UIDs] block (and everything
#    below it) shifts down by one row, opening up row 20 for the
#    (updated) [UIDs]/uidUserid/uidIssuer header and row 21 for the
#    UID/UserID/IdP sub-header, with the UID data rows following at 22-26.
$ws.Rows.Item(19).Insert()

# 3. Update the (now shifted) [UIDs] header row: uidIdP -> uidIssuer
$ws.Range("C20").Value = "uidIssuer"

# 4. Accounts table: passwords are no longer masked with a "*****"
#    formula - they now hold literal text values.
$ws.Range("C14").Value = "minderbrood"
$ws.Range("C15").Value = "joosten"
$ws.Range("C16").Value = "nolan"
$ws.Range("B17").Value = "jorrit"
$ws.Range("C17").Value = "deboer"
$ws.Range("C18").ClearContents()

# 5. Simplify the accOrg formulas (E14:E18) - drop the IF("","",...) wrapper.
$ws.Range("E14").Formula = "=`$A`$10"
$ws.Range("E15").Formula = "=`$A`$10"
$ws.Range("E16").Formula = "=`$A`$10"
$ws.Range("E17").Formula = "=`$A`$10"
$ws.Range("E18").Formula = "=`$A`$10"
